$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.395.23'
$ws.Range('E2').Value = '  +4.37%  '
$ws.Range('D3').Value = '1.723.53'
$ws.Range('E3').Value = '  +2.15%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.51'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.38%  '
$ws.Range('E6').Value = '  +0.38%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.90'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.57%  '
$ws.Range('E9').Value = '  +1.99%  '
$ws.Range('E10').Value = '  +0.89%  '
$ws.Range('E11').Value = '  +0.32%  '
$ws.Range('D12').Value = '1.969.04'
$ws.Range('E12').Value = '  +2.24%  '
$ws.Range('D13').Value = '1.730.36'
$ws.Range('E13').Value = '  +2.44%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.22'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.51%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.562'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.75%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.52'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.61%  '
$ws.Range('D17').Value = '28.346.67'
$ws.Range('E17').Value = '  +4.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '246.42'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.93%  '
$ws.Range('D19').Value = '0.0₃0749'
$ws.Range('E19').Value = '  +0.45%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.85'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.80%  '
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.59'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.05'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '149.30'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.32%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.41'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.46%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.58'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.86%  '
$ws.Range('E28').Value = '  +0.27%  '
$ws.Range('E29').Value = '  -0.24%  '
$ws.Range('E30').Value = '  +2.24%  '
$ws.Range('E31').Value = '  +2.68%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.41'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.47%  '
$ws.Range('D33').Value = '1.481.57'
$ws.Range('E33').Value = '  -4.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.23'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.44%  '
$ws.Range('E35').Value = '  -2.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.974'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.94%  '
$ws.Range('E37').Value = '  +0.57%  '
$ws.Range('E38').Value = '  -1.11%  '
$ws.Range('E39').Value = '  +1.29%  '
$ws.Range('E40').Value = '  +0.18%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '69.49'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.20%  '
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.64'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.13%  '
$ws.Range('D44').Value = '1.874.01'
$ws.Range('E44').Value = '  +1.80%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.28'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.78%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.806'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.77%  '
$ws.Range('E47').Value = '  +6.33%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '90.25'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.51%  '
$ws.Range('E49').Value = '  +2.78%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.103'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.38%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.09'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.45%  '
